$wb = $excel.ActiveWorkbook

# --- Sheet "Priorisierung" (sheet2): the only change is the active selection, ---
# --- which moved to a whole-column selection on column C.                    ---
$wsPrio = $wb.Worksheets.Item(2)
$wsPrio.Range("C1:C1048576").Select()

# --- Sheet "Tasks" (sheet3): new "erledigt am" date column, a new task row, ---
# --- and a renamed status header that moved from D1 to F2.                  ---
$wsTasks = $wb.Worksheets.Item(3)

# Preserve the old D1 text ("b ... in Bearbeitung") by writing it into F2
# *before* D1 is overwritten, so the shared string stays referenced and
# keeps its original index instead of being dropped/recreated.
$wsTasks.Range("F2").Value = "b … in Bearbeitung"

# New task row (15) - must be written before D1 is changed to "abgeschlossen"
# so that the new shared strings are appended in the same order as the target.
$wsTasks.Range("B15").Value = "für Erfolgsseiten Klasse und View Message erstellen"
$wsTasks.Range("C15").Value = "done"

# D1 becomes the new column header for the "finished" date column.
$wsTasks.Range("D1").Value = "abgeschlossen"

# Date values for the "erledigt am" column (stored as date serials, no time part)
$d0219 = (Get-Date -Year 2019 -Month 2 -Day 19).Date
$d0221 = (Get-Date -Year 2019 -Month 2 -Day 21).Date
$d0226 = (Get-Date -Year 2019 -Month 2 -Day 26).Date
$d0301 = (Get-Date -Year 2019 -Month 3 -Day 1).Date

$dateCells = @("D4", "D5", "D6", "D9", "D10", "D11", "D12", "D13", "D14", "D15")
$dateValues = @($d0219, $d0219, $d0221, $d0226, $d0226, $d0226, $d0226, $d0226, $d0226, $d0301)
for ($i = 0; $i -lt $dateCells.Length; $i++) {
    $cell = $wsTasks.Range($dateCells[$i])
    $cell.Value = $dateValues[$i]
    $cell.NumberFormat = "mm-dd-yy"
    $cell.WrapText = $true
}

# Column widths: column D narrows (it now holds short dates instead of long
# text) and a new column F is introduced for the relocated status note.
$wsTasks.Columns.Item(4).ColumnWidth = 13.333333333333334
$wsTasks.Columns.Item(6).ColumnWidth = 17.666666666666668

# Restore Tasks as the active sheet/tab and match the recorded cursor position.
$wsTasks.Activate()
$wsTasks.Range("E15").Select()
